# Update the "Create Task" bulk-upload template header row.
#
# Original layout (A1:F1): projectid | Task_role | task_mediatype | task_filedata | createdDate | modifiedDate
# New layout      (A1:I1): projectid | task_title | profile_id | Task_role | task_mediatype | task_filename | task_filepath | createdDate | modifiedDate
#
# This is accomplished the same way a human editor would do it in Excel:
#   1. Insert two new blank columns after "projectid" (new B, C) -> shifts old B..F to D..H
#   2. Insert one new blank column before the old "task_filedata" column (now F) -> shifts it (and createdDate/modifiedDate) right, new blank col becomes F
#   3. Fill in the new / changed header text
#   4. Leave selection on D8, matching the saved file

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert 2 new columns at B:C (old Task_role..modifiedDate shift from B:F to D:H)
$ws.Range("B1:C1").EntireColumn.Insert()

# Step 2: insert 1 new column at F (old task_filedata..modifiedDate shift from F:H to G:I)
$ws.Range("F1").EntireColumn.Insert()

# Step 3: populate the new columns, and rename the old "task_filedata" header (now at G1)
$ws.Range("F1").Value = "task_filename"
$ws.Range("G1").Value = "task_filepath"
$ws.Range("B1").Value = "task_title"
$ws.Range("C1").Value = "profile_id"

# Approximate the width auto-fit Excel performs for the new/changed header cells
$ws.Range("B1:C1").ColumnWidth = 8.43
$ws.Range("F1:G1").ColumnWidth = 12.75

$ws.Range("D8").Select()
